$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "harga_beli" column (D) to make
# room for the new "supplier_id" column. This shifts the old D (harga_beli)
# and E (harga_jual) columns one place to the right (to E and F).
$ws.Columns.Item(4).EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("D2").Value = "supplier_id"

# supplier_id values for each data row.
$ws.Range("D3").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 2
